# The "International Travel Ban" intervention (row 88: category
# "International Travel Ban" / parameter "Efficacy:" / R variable
# "travelban_eff") and the "School Closures" efficacy parameter
# (row 84: category "School Closures" / parameter "Efficacy:" / R
# variable "school_eff") were dropped from the model, so their rows
# are removed from the full parameter table.
#
# Delete the lower row first so the row number used for the upper
# deletion is unaffected by the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(88).Delete()
$ws.Rows.Item(84).Delete()

# Leave the selection where the author's saved file shows it: the
# last data cell of the (now shorter) table.
$ws.Range("D102").Select()
